# Updated cryptos list refresh: Price (column D) and Volume(1h) (column E).
# D-column values are plain text (e.g. "29.239.29", "1.001") that must stay
# text rather than be auto-coerced to numbers by Excel -- temporarily mark
# the cell as text ("@"), write the literal string, then restore the default
# "Normal" style so no stray number-format / style diff is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.239.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6093"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07106"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2826"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.53%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.838.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.819"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6369"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009980"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.15%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.069.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.984"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.265.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.23%  "

$ws.Range("E22").Value = "  +0.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.046"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.10%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.101"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1298"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06794"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.481"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.50%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.459"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.836"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.826"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.127"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.737"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6584"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.71%  "

$ws.Range("E37").Value = "  -0.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.232.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.50%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.761"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01765"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.592"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9316"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.31%  "

$ws.Range("E43").Value = "  +0.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.988.06"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.09%  "

$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.634"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.572"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1087"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.536"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.57%  "
